# Updates cryptos list values (price + 1h volume change) per commit
# "Updated cryptos list on Thu May 11 11:31:07 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns that now hold plain numeric-looking text (e.g. "312.92") to stay
# stored as text, matching the original inline-string cell type, instead of being
# auto-converted to numbers by Excel.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '27.483.81'
$ws.Range('E2').Value = '  -0.78%  '
$ws.Range('D3').Value = '1.825.79'
$ws.Range('E3').Value = '  -1.43%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').Value = '312.92'
$ws.Range('E5').Value = '  +0.19%  '
$ws.Range('E6').Value = '  -0.15%  '
$ws.Range('D7').Value = '0.4258'
$ws.Range('E7').Value = '  -0.49%  '
$ws.Range('E8').Value = '  +0.58%  '
$ws.Range('D9').Value = '0.07212'
$ws.Range('E9').Value = '  -1.15%  '
$ws.Range('E10').Value = '  -1.25%  '
$ws.Range('E11').Value = '  -0.86%  '
$ws.Range('D12').Value = '1.900.15'
$ws.Range('E12').Value = '  +2.48%  '
$ws.Range('D13').Value = '5.390'
$ws.Range('E13').Value = '  +1.05%  '
$ws.Range('D14').Value = '6.476'
$ws.Range('E14').Value = '  -1.10%  '
$ws.Range('D15').Value = '0.06931'
$ws.Range('E15').Value = '  -0.91%  '
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('D17').Value = '80.88'
$ws.Range('E17').Value = '  +1.68%  '
$ws.Range('D18').Value = '0.000008911'
$ws.Range('E18').Value = '  -0.38%  '
$ws.Range('E19').Value = '  -0.18%  '
$ws.Range('D20').Value = '15.38'
$ws.Range('E20').Value = '  +0.63%  '
$ws.Range('D21').Value = '28.101.95'
$ws.Range('E21').Value = '  +1.37%  '
$ws.Range('D22').Value = '5.126'
$ws.Range('E22').Value = '  +2.63%  '
$ws.Range('D23').Value = '10.86'
$ws.Range('E23').Value = '  +4.52%  '
$ws.Range('D24').Value = '2.083.76'
$ws.Range('E24').Value = '  +1.70%  '
$ws.Range('D25').Value = '1.992'
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('D26').Value = '155.37'
$ws.Range('E26').Value = '  +0.34%  '
$ws.Range('D27').Value = '18.71'
$ws.Range('E27').Value = '  +1.08%  '
$ws.Range('D28').Value = '5.147'
$ws.Range('E28').Value = '  -2.13%  '
$ws.Range('D29').Value = '114.37'
$ws.Range('E29').Value = '  -4.90%  '
$ws.Range('E30').Value = '  -4.82%  '
$ws.Range('D31').Value = '0.08900'
$ws.Range('E31').Value = '  -0.19%  '
$ws.Range('D32').Value = '0.7491'
$ws.Range('E32').Value = '  -1.08%  '
$ws.Range('B33').Value = 'HuobiToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D33').Value = '2.976'
$ws.Range('E33').Value = '  +0.19%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '4.544'
$ws.Range('E34').Value = '  +0.72%  '
$ws.Range('D35').Value = '1.120'
$ws.Range('E35').Value = '  -0.22%  '
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('D37').Value = '1.085'
$ws.Range('E37').Value = '  -2.06%  '
$ws.Range('D38').Value = '0.05269'
$ws.Range('E38').Value = '  -2.96%  '
$ws.Range('E39').Value = '  -0.17%  '
$ws.Range('D40').Value = '2.793'
$ws.Range('E40').Value = '  -1.17%  '
$ws.Range('D41').Value = '0.5077'
$ws.Range('E41').Value = '  -0.18%  '
$ws.Range('D42').Value = '0.1655'
$ws.Range('E42').Value = '  -1.05%  '
$ws.Range('D43').Value = '6.443'
$ws.Range('E43').Value = '  -2.76%  '
$ws.Range('D44').Value = '8.349'
$ws.Range('E44').Value = '  -0.68%  '
$ws.Range('D45').Value = '10.49'
$ws.Range('E45').Value = '  +1.34%  '
$ws.Range('D46').Value = '106.50'
$ws.Range('E46').Value = '  +0.33%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '0.06463'
$ws.Range('E47').Value = '  -0.97%  '
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').Value = '0.4676'
$ws.Range('E48').Value = '  +0.15%  '
$ws.Range('E49').Value = '  -0.18%  '
$ws.Range('E50').Value = '  -0.27%  '
$ws.Range('D51').Value = '64.06'
$ws.Range('E51').Value = '  -0.51%  '
